$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $savedStyle
}

Set-TextValue $ws.Range("D2") '24.529.36'
$ws.Range("E2").Value = '  +3.51%  '
Set-TextValue $ws.Range("D3") '1.692.74'
$ws.Range("E3").Value = '  +2.15%  '
Set-TextValue $ws.Range("D4") '1.004'
$ws.Range("E4").Value = '  +0.24%  '
Set-TextValue $ws.Range("D5") '316.36'
$ws.Range("E5").Value = '  +2.20%  '
$ws.Range("E6").Value = '  +0.01%  '
Set-TextValue $ws.Range("D7") '0.3944'
$ws.Range("E7").Value = '  +1.46%  '
Set-TextValue $ws.Range("D8") '0.4014'
$ws.Range("E8").Value = '  +1.59%  '
Set-TextValue $ws.Range("D9") '1.532'
$ws.Range("E9").Value = '  +6.80%  '
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D10") '54.43'
$ws.Range("E10").Value = '  +10.78%  '
$ws.Range("B11").Value = 'BinanceUSD'
$ws.Range("C11").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range("D11") '1.004'
$ws.Range("E11").Value = '  +0.15%  '
Set-TextValue $ws.Range("D12") '0.08752'
$ws.Range("E12").Value = '  +1.26%  '
Set-TextValue $ws.Range("D13") '7.230'
$ws.Range("E13").Value = '  +8.39%  '
Set-TextValue $ws.Range("D14") '23.23'
$ws.Range("E14").Value = '  +2.87%  '
Set-TextValue $ws.Range("D15") '0.00001315'
$ws.Range("E15").Value = '  +0.02%  '
Set-TextValue $ws.Range("D16") '7.591'
$ws.Range("E16").Value = '  +5.04%  '
Set-TextValue $ws.Range("D17") '1.694.85'
$ws.Range("E17").Value = '  +2.30%  '
Set-TextValue $ws.Range("D18") '101.05'
$ws.Range("E18").Value = '  +1.41%  '
Set-TextValue $ws.Range("D19") '0.07010'
$ws.Range("E19").Value = '  +3.56%  '
Set-TextValue $ws.Range("D20") '19.62'
Set-TextValue $ws.Range("D21") '6.854'
$ws.Range("E21").Value = '  +2.94%  '
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("E23").Value = '  +1.01%  '
Set-TextValue $ws.Range("D24") '24.525.58'
$ws.Range("E24").Value = '  +3.55%  '
Set-TextValue $ws.Range("D25") '3.044'
$ws.Range("E25").Value = '  +7.89%  '
Set-TextValue $ws.Range("D26") '2.323'
$ws.Range("E26").Value = '  +0.14%  '
Set-TextValue $ws.Range("D27") '22.37'
$ws.Range("E27").Value = '  +2.99%  '
Set-TextValue $ws.Range("D28") '159.46'
$ws.Range("E28").Value = '  +0.65%  '
Set-TextValue $ws.Range("D29") '5.190'
$ws.Range("E29").Value = '  +1.34%  '
Set-TextValue $ws.Range("D30") '134.23'
$ws.Range("E30").Value = '  +3.69%  '
Set-TextValue $ws.Range("D31") '7.484'
$ws.Range("E31").Value = '  +16.35%  '
Set-TextValue $ws.Range("D32") '1.883.10'
$ws.Range("E32").Value = '  +2.38%  '
Set-TextValue $ws.Range("D33") '1.088'
$ws.Range("E33").Value = '  -4.01%  '
Set-TextValue $ws.Range("D34") '7.330'
$ws.Range("E34").Value = '  +11.95%  '
Set-TextValue $ws.Range("D35") '0.08516'
$ws.Range("E35").Value = '  -0.75%  '
Set-TextValue $ws.Range("D36") '11.41'
$ws.Range("E36").Value = '  +10.01%  '
Set-TextValue $ws.Range("D37") '1.979'
$ws.Range("E37").Value = '  -0.58%  '
Set-TextValue $ws.Range("D38") '0.2725'
$ws.Range("E38").Value = '  +3.03%  '
Set-TextValue $ws.Range("D39") '14.54'
$ws.Range("E39").Value = '  +1.23%  '
Set-TextValue $ws.Range("D40") '0.02750'
$ws.Range("E40").Value = '  +9.05%  '
Set-TextValue $ws.Range("D41") '0.09005'
$ws.Range("E41").Value = '  +2.66%  '
Set-TextValue $ws.Range("D42") '1.461'
$ws.Range("E42").Value = '  +1.11%  '
Set-TextValue $ws.Range("D43") '0.7665'
$ws.Range("E43").Value = '  +1.67%  '
Set-TextValue $ws.Range("D44") '0.7183'
$ws.Range("E44").Value = '  +2.49%  '
Set-TextValue $ws.Range("D45") '15.32'
$ws.Range("E45").Value = '  +3.20%  '
Set-TextValue $ws.Range("D46") '2.512'
$ws.Range("E46").Value = '  +4.41%  '
Set-TextValue $ws.Range("D47") '4.221'
$ws.Range("E47").Value = '  +3.01%  '
Set-TextValue $ws.Range("D48") '1.002'
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("E49").Value = '  +13.47%  '
Set-TextValue $ws.Range("D50") '140.91'
$ws.Range("E50").Value = '  +2.15%  '
Set-TextValue $ws.Range("D51") '0.08013'
$ws.Range("E51").Value = '  +3.38%  '
